$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data order: sorted descending by value, with Swedish and Uzbek removed
$data = @(
    @("Chinese", 20.69999591507272),
    @("English", 20.46865504318676),
    @("Spanish", 6.007631145698993),
    @("Arabic", 4.192068016298176),
    @("German", 3.868047895741215),
    @("Japanese", 3.540727294876807),
    @("Malay-Indonesian", 3.342566871720621),
    @("Russian", 2.873835974480956),
    @("Portuguese", 2.692229765739679),
    @("French", 2.413157834688444),
    @("Turkish", 2.057318389698041),
    @("Italian", 1.759333203911673),
    @("Korean", 1.654330262287356),
    @("Dutch", 1.145607077290907),
    @("Polish", 0.9772779184520434),
    @("Persian", 0.9507849578428667),
    @("Urdu", 0.9233230964212991),
    @("Thai", 0.9187433561642289),
    @("Bengali", 0.9146009165946761),
    @("Vietnamese", 0.8931452039523774)
)

$row = 2
foreach ($pair in $data) {
    $ws.Cells.Item($row, 1).Value = $pair[0]
    $ws.Cells.Item($row, 2).Value = $pair[1]
    $row++
}

# Remove the two now-unused trailing rows (previously rows 22 and 23)
$ws.Range("A22:B23").Clear()
